$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the D/E columns we touch stay plain text (matches the original
# inlineStr cells) instead of Excel auto-coercing numeric-looking strings
# like "1.00" into numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "34.048.55"
$ws.Range("E2").Value = "  +10.42%  "

$ws.Range("D3").Value = "1.811.22"
$ws.Range("E3").Value = "  +7.83%  "

$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").Value = "228.72"
$ws.Range("E5").Value = "  +4.21%  "

$ws.Range("D6").Value = "0.571"
$ws.Range("E6").Value = "  +6.40%  "

$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.43%  "

$ws.Range("D8").Value = "31.77"
$ws.Range("E8").Value = "  +9.70%  "

$ws.Range("D9").Value = "46.59"
$ws.Range("E9").Value = "  +5.55%  "

$ws.Range("D10").Value = "0.284"
$ws.Range("E10").Value = "  +7.92%  "

$ws.Range("D11").Value = "0.0670"
$ws.Range("E11").Value = "  +3.89%  "

$ws.Range("E12").Value = "  +2.95%  "

$ws.Range("D13").Value = "2.081.60"
$ws.Range("E13").Value = "  +8.42%  "

$ws.Range("D14").Value = "1.806.90"
$ws.Range("E14").Value = "  +7.21%  "

$ws.Range("D15").Value = "0.641"
$ws.Range("E15").Value = "  +6.13%  "

$ws.Range("D16").Value = "34.065.01"
$ws.Range("E16").Value = "  +10.60%  "

$ws.Range("D17").Value = "10.10"
$ws.Range("E17").Value = "  -0.47%  "

$ws.Range("D18").Value = "4.24"
$ws.Range("E18").Value = "  +4.94%  "

$ws.Range("D19").Value = "69.93"
$ws.Range("E19").Value = "  +5.85%  "

$ws.Range("D20").Value = "257.03"
$ws.Range("E20").Value = "  +5.58%  "

$ws.Range("D21").Value = "0.0₃0752"
$ws.Range("E21").Value = "  +4.33%  "

$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").Value = "10.56"
$ws.Range("E23").Value = "  +6.01%  "

$ws.Range("D24").Value = "4.30"
$ws.Range("E24").Value = "  +1.47%  "

$ws.Range("D25").Value = "2.21"
$ws.Range("E25").Value = "  +1.98%  "

$ws.Range("D26").Value = "159.74"
$ws.Range("E26").Value = "  +0.35%  "

$ws.Range("D27").Value = "16.63"
$ws.Range("E27").Value = "  +5.04%  "

$ws.Range("D28").Value = "0.117"
$ws.Range("E28").Value = "  +4.15%  "

$ws.Range("D29").Value = "7.02"
$ws.Range("E29").Value = "  +5.05%  "

$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.39%  "

$ws.Range("D31").Value = "3.86"
$ws.Range("E31").Value = "  +11.44%  "

$ws.Range("D32").Value = "0.0523"
$ws.Range("E32").Value = "  +6.00%  "

$ws.Range("D33").Value = "1.21"
$ws.Range("E33").Value = "  +6.09%  "

$ws.Range("D34").Value = "3.58"
$ws.Range("E34").Value = "  +8.38%  "

$ws.Range("D35").Value = "1.512.40"
$ws.Range("E35").Value = "  -0.57%  "

$ws.Range("D36").Value = "1.78"
$ws.Range("E36").Value = "  +1.22%  "

$ws.Range("E37").Value = "  +5.19%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.0189"
$ws.Range("E38").Value = "  +5.80%  "

$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "0.628"
$ws.Range("E39").Value = "  +4.36%  "

$ws.Range("D40").Value = "82.86"
$ws.Range("E40").Value = "  -1.39%  "

$ws.Range("D41").Value = "2.77"
$ws.Range("E41").Value = "  +4.47%  "

$ws.Range("D42").Value = "2.38"
$ws.Range("E42").Value = "  +3.88%  "

$ws.Range("D43").Value = "0.902"
$ws.Range("E43").Value = "  +7.53%  "

$ws.Range("D44").Value = "2.10"
$ws.Range("E44").Value = "  +3.53%  "

$ws.Range("D45").Value = "0.0520"
$ws.Range("E45").Value = "  +4.04%  "

$ws.Range("E46").Value = "  +5.40%  "

$ws.Range("D47").Value = "1.968.80"
$ws.Range("E47").Value = "  +8.53%  "

$ws.Range("D48").Value = "5.84"
$ws.Range("E48").Value = "  +5.38%  "

$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.31%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "11.98"
$ws.Range("E50").Value = "  +14.30%  "

$ws.Range("D51").Value = "51.42"
$ws.Range("E51").Value = "  +1.41%  "
